$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F2").Value = 4748  # was 4746
$wsExhibition.Range("F6").Value = 569  # was 567
$wsExhibition.Range("F7").Value = 572  # was 570
$wsExhibition.Range("F8").Value = 429  # was 426
$wsExhibition.Range("F9").Value = 150  # was 148
$wsExhibition.Range("F10").Value = 1819  # was 1816
$wsExhibition.Range("F11").Value = 1416  # was 1411
$wsExhibition.Range("F13").Value = 1681  # was 1677
$wsExhibition.Range("F14").Value = 29  # was 28
$wsExhibition.Range("F16").Value = 636  # was 634
$wsExhibition.Range("F17").Value = 23  # was 22
$wsExhibition.Range("F18").Value = 51  # was 50
$wsExhibition.Range("F21").Value = 66  # was 65
$wsExhibition.Range("F23").Value = 16  # was 15
$wsExhibition.Range("F25").Value = 54  # was 53
$wsExhibition.Range("F27").Value = 4259  # was 4232
$wsExhibition.Range("F29").Value = 791  # was 789
$wsExhibition.Range("F31").Value = 2082  # was 2051
$wsExhibition.Range("F32").Value = 62  # was 61
$wsExhibition.Range("F33").Value = 1970  # was 1963

# Sheet 2: 演出 (Show)
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F3").Value = 63  # was 62

# Sheet 4: 全部类型 (All types)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 4748  # was 4746
$wsAll.Range("F6").Value = 569  # was 567
$wsAll.Range("F7").Value = 572  # was 570
$wsAll.Range("F9").Value = 429  # was 426
$wsAll.Range("F10").Value = 150  # was 148
$wsAll.Range("F11").Value = 1819  # was 1816
$wsAll.Range("F12").Value = 1416  # was 1411
$wsAll.Range("F14").Value = 1681  # was 1677
$wsAll.Range("F15").Value = 29  # was 28
$wsAll.Range("F17").Value = 636  # was 634
$wsAll.Range("F18").Value = 23  # was 22
$wsAll.Range("F19").Value = 51  # was 50
$wsAll.Range("F22").Value = 66  # was 65
$wsAll.Range("F24").Value = 16  # was 15
$wsAll.Range("F26").Value = 54  # was 53
$wsAll.Range("F28").Value = 4259  # was 4232
$wsAll.Range("F29").Value = 63  # was 62
$wsAll.Range("F32").Value = 791  # was 789
$wsAll.Range("F34").Value = 2082  # was 2051
$wsAll.Range("F35").Value = 62  # was 61
$wsAll.Range("F36").Value = 1970  # was 1963
